$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("key-to-input-workbooks")

# Clear existing contents (sheet previously had rows 1-30; new table spans 1-39)
$ws.Cells.Clear()

$colA = @{}
$colB = @{}
$colA[1] = "L-curve ID"
$colB[1] = "Filename"
$colA[2] = 1
$colB[2] = "15-genes_28-edges_db5_Sigmoid_estimation_no-missing-values_L-curve.xlsx SEA120-14"
$colA[3] = 2
$colB[3] = "15-genes_28-edges_db5_Sigmoid_estimation_missing-values_L-curve.xlsx SEA120-14"
$colA[4] = 3
$colB[4] = "15-genes_28-edges_db5_Sigmoid_estimation_no-missing-values_L-curve.xlsx SEA120-15"
$colA[5] = 4
$colB[5] = "15-genes_28-edges_db5_Sigmoid_estimation_missing-values_L-curve.xlsx SEA120-15"
$colA[6] = 5
$colB[6] = "15-genes_28-edges_db5_Sigmoid_estimation_missing-values_fixb.xlsx"
$colA[7] = 6
$colB[7] = "15-genes_28-edges_db5_Sigmoid_estimation_missing-values_fixP.xlsx"
$colA[8] = 7
$colB[8] = "15-genes_28-edges_db5_Sigmoid_estimation_missing-values_fixP-fixb.xlsx"
$colA[9] = 8
$colB[9] = "15-genes_28-edges_db5_Sigmoid_estimation_missing-values_L-curve_no-GLN3-data.xlsx"
$colA[10] = 9
$colB[10] = "15-genes_28-edges_db5_Sigmoid_estimation_missing-values_L-curve_no-GLN3-ZAP1-data.xlsx"
$colA[11] = 10
$colB[11] = "15-genes_28-edges_db5_Sigmoid_estimation_missing-values_L-curve_no-GLN3-ZAP1-HAP4-data.xlsx"
$colA[12] = 11
$colB[12] = "15-genes_28-edges_db5_Sigmoid_estimation_missing-values_L-curve_no-GLN3-ZAP1-HAP4-HMO1-CIN5-data.xlsx"
$colA[13] = 12
$colB[13] = "15-genes_28-edges_db5_Sigmoid_estimation_missing-values_L-curve_no-GLN3-ZAP1-HAP4-HMO1-data.xlsx"
$colA[14] = 13
$colB[14] = "14-genes_27-edges_db5_Sigmoid_estimation_missing-values_delGLN3.xlsx"
$colA[15] = 14
$colB[15] = "13-genes_26-edges_db5_Sigmoid_estimation_missing-values_delGLN3-delZAP1.xlsx"
$colA[16] = 15
$colB[16] = "12-genes_25-edges_db5_Sigmoid_estimation_missing-values_delGLN3--delZAP1-delGCR2.xlsx"
$colA[17] = 16
$colB[17] = "11-genes_24-edges_db5_Sigmoid_estimation_missing-values_delGLN3-delZAP1-delGCR2-delACE2.xlsx"
$colA[18] = 17
$colB[18] = "10-genes_22-edges_db5_Sigmoid_estimation_missing-values_delGLN3-delZAP1-delGCR2-delACE2-delSWI5.xlsx"
$colA[19] = 18
$colB[19] = "9-genes_20-edges_db5_Sigmoid_estimation_missing-values_delGLN3-delZAP1-delGCR2-delACE2-delSWI5-delASH1.xlsx"
$colA[20] = 19
$colB[20] = "8-genes_17-edges_db5_Sigmoid_estimation_missing-values_delGLN3-delZAP1-delGCR2-delACE2-delSWI5-delASH1-delYOX1.xlsx"
$colA[21] = 20
$colB[21] = "7-genes_14-edges_db5_Sigmoid_estimation_missing-values_delGLN3-delZAP1-delGCR2-delACE2-delSWI5-delASH1-delYOX1-delYHP1.xlsx"
$colA[22] = 21
$colB[22] = "6-genes_11-edges_db5_Sigmoid_estimation_missing-values_delGLN3-delZAP1-delGCR2-delACE2-delSWI5-delASH1-delYOX1-delYHP1-delSFP1.xlsx"
$colA[23] = 22
$colB[23] = "5-genes_9-edges_db5_Sigmoid_estimation_missing-values_dGLN3-dZAP1-dGCR2-dACE2-dSWI5-dASH1-dYOX1-dYHP1-dSFP1-dSWI4.xlsx"
$colA[24] = 23
$colB[24] = "4-genes_7-edges_db5_Sigmoid_estimation_missing-values_dGLN3-dZAP1-dGCR2-dACE2-dSWI5-dASH1-dYOX1-dYHP1-dSFP1-dSWI4-dSTB5.xlsx"
$colA[25] = 24
$colB[25] = "3-genes_4-edges_db5_Sigmoid_estimation_missing-values_dGLN3-dZAP1-dGCR2-dACE2-dSWI5-dASH1-dYOX1-dYHP1-dSFP1-dSWI4-dSTB5-dMSN2.xlsx"
$colA[26] = 25
$colB[26] = "15-genes_28-edges_db5-random2_Sigmoid_estimation_missing-values_L-curve.xlsx"
$colA[27] = 26
$colB[27] = "15-genes_28-edges_db5-random3_Sigmoid_estimation_missing-values_L-curve.xlsx"
$colA[28] = 27
$colB[28] = "15-genes_28-edges_db5-random7-fam_Sigmoid_estimation_missing-values_L-curve.xlsx"
$colA[29] = 28
$colB[29] = "15-genes_28-edges_db5-random9_Sigmoid_estimation_missing-values_L-curve.xlsx"
$colA[30] = 29
$colB[30] = "15-genes_28-edges_db5-random12-fam_Sigmoid_estimation_missing-values_L-curve.xlsx"
$colA[31] = 30
$colB[31] = "15-genes_28-edges_db5-random15_Sigmoid_estimation_missing-values_L-curve.xlsx"
$colA[32] = 31
$colB[32] = "15-genes_28-edges_db5-random16-fam_Sigmoid_estimation_missing-values_L-curve.xlsx"
$colA[33] = 32
$colB[33] = "15-genes_28-edges_db5-random24-fam_Sigmoid_estimation_missing-values_L-curve.xlsx"
$colA[34] = 33
$colB[34] = "15-genes_28-edges_db5-random31_Sigmoid_estimation_missing-values_L-curve.xlsx"
$colA[35] = 34
$colB[35] = "16-genes_36-edges_db1_Sigmoid_estimation_missing-values_L-curve.xlsx"
$colA[36] = 35
$colB[36] = "14-gene_25-edges_db2_Sigmoid_estimation_missing-values_L-curve.xlsx"
$colA[37] = 36
$colB[37] = "17-genes_32-edges_db3_Sigmoid_estimation_missing-values_L-curve.xlsx"
$colA[38] = 37
$colB[38] = "14-genes_35-edges_db4_Sigmoid_estimation_missing-values_L-curve.xlsx"
$colA[39] = 38
$colB[39] = "16-genes_27-edges_db6_Sigmoid_estimation_missing-values_L-curve.xlsx"

$writeOrder = @(6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21, 22, 23, 24, 25, 26, 27, 28, 29, 30, 31, 32, 33, 34, 35, 36, 37, 38, 39, 1, 2, 3, 4, 5)
foreach ($r in $writeOrder) {
    $ws.Cells.Item($r, 1).Value = $colA[$r]
    $ws.Cells.Item($r, 2).Value = $colB[$r]
}

# Restore the view to the top-left cell (the prior explicit selection on A2
# is no longer meaningful once the table has been rebuilt).
$ws.Activate() | Out-Null
$ws.Range("A1").Select() | Out-Null